$d = $word.ActiveDocument

# --- Change 1: merge the split "TUE Sep 18" / " 11:30:25 IST 2018" runs into
# a single run, same as retyping the identical text over the found range. ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("TUE Sep 18 11:30:25 IST 2018", $true, $false, $false, $false, $false, `
              $true, 1, $false, "TUE Sep 18 11:30:25 IST 2018", 2) | Out-Null

# --- Change 2: append the new 27/09/2018 purchase-details block (MAMATHA
# CHICK IN) after the last entry in the document ("Amount Received mode -
# CASH"), before the trailing blank paragraphs. ---

$lastEntry = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text.TrimEnd([char]13)
    if ($t -like "Amount Received mode*- CASH") {
        $lastEntry = $para
    }
}
if ($lastEntry -eq $null) {
    throw "Could not locate the 'Amount Received mode ... - CASH' paragraph to anchor the new block."
}

$r = $lastEntry.Range
$insertionPoint = $d.Range($r.End - 1, $r.End - 1)

$runFonts = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>'

$newParasXml = (
  '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $runFonts + '<w:b/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $runFonts + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:t>TUE Sep 25</w:t></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:t xml:space="preserve"> 11:41:46 IST 2018</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $runFonts + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:t>Person Name</w:t></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:tab/></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:tab/></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:tab/></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:tab/><w:t>- NAGAMMA</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $runFonts + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:t>---------------------------------------------------------------</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $runFonts + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:t>Item Name</w:t></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:tab/></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:tab/></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:tab/></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:tab/><w:t>- BEET</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $runFonts + '<w:color w:val="FF0000"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $runFonts + '<w:color w:val="FF0000"/></w:rPr><w:t>Amount Received</w:t></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '<w:color w:val="FF0000"/></w:rPr><w:tab/></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '<w:color w:val="FF0000"/></w:rPr><w:tab/></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '<w:color w:val="FF0000"/></w:rPr><w:tab/><w:t>- 550</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $runFonts + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:t>Amount Received mode</w:t></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:tab/></w:r>' +
    '<w:r><w:rPr>' + $runFonts + '</w:rPr><w:tab/><w:t>- CASH AND CLEAD</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $runFonts + '</w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $runFonts + '<w:b/></w:rPr></w:pPr></w:p>'
)

$packageXml = (
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" ' +
      'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' + $newParasXml + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
)

$insertionPoint.InsertXML($packageXml) | Out-Null

Write-Output ("Paragraph count after edit: " + $d.Paragraphs.Count)
